$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.673.00"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.638.74"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'212.39"
$ws.Range("D6").Value = "'0.523"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'23.09"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "1.871.65"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "1.636.23"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Value = "'4.06"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  -5.52%  "
$ws.Range("D16").Value = "'64.66"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "27.664.26"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "'230.57"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "'10.22"
$ws.Range("E23").Value = "  +4.20%  "
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").Value = "'151.43"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "1.457.83"
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "'0.566"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").Value = "'0.878"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'0.895"
$ws.Range("E40").Value = "  +9.14%  "
$ws.Range("D41").Value = "'69.79"
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "1.781.23"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  +3.34%  "
$ws.Range("D49").Value = "'86.89"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  -0.16%  "
